# Applies the "Updated symbol list" GitHub Actions refresh (Sat Dec 17 15:04:26 UTC 2022)
# to cryptos.xlsx: refreshed Price (D) / Hora (G) columns for every listed coin, and two
# rows (42/43) swap contents because BKEXToken/CEJI traded ranking places upstream.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every data value as text (inline/shared string), including
# cells that look numeric ("238.65", "15", ...). Assigning such a literal straight to
# Range.Value lets Excel auto-convert it to a real number, which would lose formatting
# (e.g. trailing zeros in "0.05630"). Routing the literal through a `="..."` text formula
# and then collapsing it to a plain value via Copy/PasteSpecial keeps it a text cell.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Formula = '="' + $text + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

Set-TextValue 'D2' '238.65'
Set-TextValue 'G2' '15'
Set-TextValue 'D3' '21.78'
Set-TextValue 'G3' '15'
Set-TextValue 'D4' '3.904'
Set-TextValue 'G4' '15'
Set-TextValue 'D5' '5.388'
Set-TextValue 'G5' '15'
Set-TextValue 'D6' '0.05630'
Set-TextValue 'G6' '15'
Set-TextValue 'D7' '6.477'
Set-TextValue 'G7' '15'
Set-TextValue 'D8' '3.342'
Set-TextValue 'G8' '15'
Set-TextValue 'D9' '0.7966'
Set-TextValue 'G9' '15'
Set-TextValue 'G10' '15'
Set-TextValue 'D11' '0.01171'
Set-TextValue 'G11' '15'
Set-TextValue 'G12' '15'
Set-TextValue 'D13' '0.07341'
Set-TextValue 'G13' '15'
Set-TextValue 'D14' '0.03145'
Set-TextValue 'G14' '15'
Set-TextValue 'D15' '0.02980'
Set-TextValue 'G15' '15'
Set-TextValue 'D16' '0.09233'
Set-TextValue 'G16' '15'
Set-TextValue 'D17' '0.001672'
Set-TextValue 'G17' '15'
Set-TextValue 'D18' '3.254'
Set-TextValue 'G18' '15'
Set-TextValue 'D19' '0.04773'
Set-TextValue 'G19' '15'
Set-TextValue 'D20' '0.006246'
Set-TextValue 'G20' '15'
Set-TextValue 'D21' '0.005075'
Set-TextValue 'G21' '15'
Set-TextValue 'D22' '0.001052'
Set-TextValue 'G22' '15'
Set-TextValue 'D23' '0.0001502'
Set-TextValue 'G23' '15'
Set-TextValue 'D24' '0.0004005'
Set-TextValue 'G24' '15'
Set-TextValue 'G25' '15'
Set-TextValue 'G26' '15'
Set-TextValue 'D27' '0.1053'
Set-TextValue 'G27' '15'
Set-TextValue 'G28' '15'
Set-TextValue 'G29' '15'
Set-TextValue 'G30' '15'
Set-TextValue 'G31' '15'
Set-TextValue 'G32' '15'
Set-TextValue 'G33' '15'
Set-TextValue 'G34' '15'
Set-TextValue 'G35' '15'
Set-TextValue 'G36' '15'
Set-TextValue 'G37' '15'
Set-TextValue 'G38' '15'
Set-TextValue 'G39' '15'
Set-TextValue 'D40' '0.04077'
Set-TextValue 'G40' '15'
Set-TextValue 'D41' '0.006941'
Set-TextValue 'G41' '15'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 'D42' '0.1038'
$ws.Range('E42').Value = '41BKEXTokenBKK'
Set-TextValue 'G42' '15'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue 'D43' '0.003161'
$ws.Range('E43').Value = '42CEJICEJI'
Set-TextValue 'G43' '15'
Set-TextValue 'D44' '0.009111'
Set-TextValue 'G44' '15'
Set-TextValue 'D45' '0.00005448'
Set-TextValue 'G45' '15'
Set-TextValue 'D46' '0.00000000751'
Set-TextValue 'G46' '15'
Set-TextValue 'D47' '0.6760'
Set-TextValue 'G47' '15'
Set-TextValue 'D48' '0.03753'
Set-TextValue 'G48' '15'
Set-TextValue 'D49' '0.00002102'
Set-TextValue 'G49' '15'
Set-TextValue 'D50' '0.01011'
Set-TextValue 'G50' '15'
Set-TextValue 'G51' '15'

$excel.CutCopyMode = $false
Write-Host "Applied cryptos.xlsx symbol-list update (Sat Dec 17 15:04:26 UTC 2022)"
